# Reset the DataFrame index written into column A:
# - Remove the old "index" header cell at A1 (header row now starts at B1)
# - Replace the string labels (BHP-FY-2003 ... BHP-FY-2022) in A2:A21 with a
#   plain 0-based numeric index (0 .. 19), keeping the existing cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header text in A1 completely (content + formatting) so no
# cell remains there at all.
$ws.Range("A1").Clear()

# Replace the string index values in A2:A21 with integers 0..19.
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}
